# Generate Report for Handback
#
# Populates the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns on each locale sheet with the freshly
# generated handback artifacts, flips the Overview status from
# "Ready for handoff" to "Handed back: in sync with en-US", and widens a
# few columns so the longer handback file names remain readable.

$wb = $excel.ActiveWorkbook

$srcMdUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6d992497c3ca85c678b12b62ef2857972c7b2f57/e2e/5d24a5f0-095a-4f6c-ba13-767e5b80d782.md"
$srcMdUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6d992497c3ca85c678b12b62ef2857972c7b2f57/e2e/84cb22ab-5658-4cb6-b7d7-fbe251bee46d.md"
$srcMdName1 = "5d24a5f0-095a-4f6c-ba13-767e5b80d782.md"
$srcMdName2 = "84cb22ab-5658-4cb6-b7d7-fbe251bee46d.md"

$wideColumnWidth = 29.166666666666668
$maxColumnWidth  = 39.166666666666664

# ---------------------------------------------------------------------------
# Overview sheet: the handback is complete and in sync, update the status
# shown for both locales.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
$wsOverview.Columns.Item(5).ColumnWidth = $wideColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $wideColumnWidth

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $srcMdUrl1, "", "", $srcMdName1)
$wsZh.Range("J2").Value = "5d24a5f0-095a-4f6c-ba13-767e5b80d782.b0c7a6fa9b68ec5f8a10893d370c0d2f11b82a62.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-31 19:11:36"

$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $srcMdUrl2, "", "", $srcMdName2)
$wsZh.Range("J3").Value = "84cb22ab-5658-4cb6-b7d7-fbe251bee46d.313c14a06a567d9a49dcb727d9d10f26a6f8a805.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-31 19:11:36"

$wsZh.Columns.Item(3).ColumnWidth = $wideColumnWidth
$wsZh.Columns.Item(9).ColumnWidth = $maxColumnWidth
$wsZh.Columns.Item(10).ColumnWidth = $maxColumnWidth

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $srcMdUrl1, "", "", $srcMdName1)
$wsDe.Range("J2").Value = "5d24a5f0-095a-4f6c-ba13-767e5b80d782.b0c7a6fa9b68ec5f8a10893d370c0d2f11b82a62.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-31 19:11:44"

$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $srcMdUrl2, "", "", $srcMdName2)
$wsDe.Range("J3").Value = "84cb22ab-5658-4cb6-b7d7-fbe251bee46d.313c14a06a567d9a49dcb727d9d10f26a6f8a805.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-31 19:11:44"

$wsDe.Columns.Item(3).ColumnWidth = $wideColumnWidth
$wsDe.Columns.Item(9).ColumnWidth = $maxColumnWidth
$wsDe.Columns.Item(10).ColumnWidth = $maxColumnWidth
